# ============================================================================
# Scheduled-runner market data refresh for the Leve-profit workbook.
#
# For every (sheet, row) below, columns H:N hold point-in-time Market Board
# data (current average prices) and the derived Leve profit figures. This
# script overwrites those columns with the freshly-pulled values; columns
# A:G (leve/item metadata) are untouched. A few rows lose their NQ-profit
# figure entirely (column M) because NQ is no longer sellable for that item -
# those cells are cleared instead of zeroed so the row keeps "no data" rather
# than a misleading 0.
# ============================================================================

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# row 19
$ws.Range("H19").Value = 704.2222
$ws.Range("I19").Value = 577.6667
$ws.Range("J19").Value = 767.5
$ws.Range("K19").Value = 577.6667
$ws.Range("L19").Value = 767.5
$ws.Range("M19").Value = -402.6667
$ws.Range("N19").Value = -1117.5
# row 32
$ws.Range("H32").Value = 4454.8
$ws.Range("I32").Value = 5866
$ws.Range("J32").Value = 4205.7646
$ws.Range("K32").Value = 5866
$ws.Range("L32").Value = 4205.7646
$ws.Range("M32").Value = -5540
$ws.Range("N32").Value = -4857.7646
# row 96
$ws.Range("H96").Value = 1725.1818
$ws.Range("I96").Value = 387.42856
$ws.Range("J96").Value = 4066.25
$ws.Range("K96").Value = 1162.28568
$ws.Range("L96").Value = 12198.75
$ws.Range("M96").Value = 210.71432
$ws.Range("N96").Value = -14944.75
# row 116
$ws.Range("H116").Value = 37208396
$ws.Range("I116").Value = 28515006
$ws.Range("J116").Value = 55561104
$ws.Range("K116").Value = 28515006
$ws.Range("L116").Value = 55561104
$ws.Range("M116").Value = -28511564
$ws.Range("N116").Value = -55567988

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# row 32
$ws.Range("H32").Value = 4849.15
$ws.Range("I32").Value = 2867.5356
$ws.Range("J32").Value = 9472.916999999999
$ws.Range("K32").Value = 2867.5356
$ws.Range("L32").Value = 9472.916999999999
$ws.Range("M32").Value = -2580.5356
$ws.Range("N32").Value = -10046.917
# row 61
$ws.Range("H61").Value = 23011.715
$ws.Range("I61").Value = 18612.25
$ws.Range("J61").Value = 28877.666
$ws.Range("K61").Value = 18612.25
$ws.Range("L61").Value = 28877.666
$ws.Range("M61").Value = -18400.25
$ws.Range("N61").Value = -29301.666
# row 63
$ws.Range("H63").Value = 3248.889
$ws.Range("J63").Value = 3228
$ws.Range("L63").Value = 3228
$ws.Range("N63").Value = -4600
# row 66
$ws.Range("H66").Value = 3248.889
$ws.Range("J66").Value = 3228
$ws.Range("L66").Value = 16140
$ws.Range("N66").Value = -23004
# row 102
$ws.Range("H102").Value = 596709
$ws.Range("J102").Value = 1157
$ws.Range("L102").Value = 1157
$ws.Range("N102").Value = -4401
# row 122
$ws.Range("H122").Value = 4949.6875
$ws.Range("I122").Value = 2885.9092
$ws.Range("K122").Value = 8657.7276
$ws.Range("M122").Value = -6207.7276
# row 132
$ws.Range("H132").Value = 15950.282
$ws.Range("I132").Value = 16962.639
$ws.Range("K132").Value = 50887.917
$ws.Range("M132").Value = -48357.917
# row 136
$ws.Range("H136").Value = 23011.715
$ws.Range("I136").Value = 18612.25
$ws.Range("J136").Value = 28877.666
$ws.Range("K136").Value = 55836.75
$ws.Range("L136").Value = 86632.99800000001
$ws.Range("M136").Value = -53286.75
$ws.Range("N136").Value = -91732.99800000001
# row 140
$ws.Range("H140").Value = 114433.336
$ws.Range("J140").Value = 114433.336
$ws.Range("L140").Value = 114433.336
$ws.Range("N140").Value = -124793.336

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# row 140
$ws.Range("H140").Value = 98057.836
$ws.Range("J140").Value = 98057.836
$ws.Range("L140").Value = 98057.836
$ws.Range("N140").Value = -108417.836

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# row 11
$ws.Range("H11").Value = 1621.8572
$ws.Range("I11").Value = 101
$ws.Range("K11").Value = 101
$ws.Range("M11").Value = 39
# row 16
$ws.Range("H16").Value = 2023.6666
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 2023.6666
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 2023.6666
$ws.Range("M16").ClearContents()
$ws.Range("N16").Value = -2597.6666
# row 31
$ws.Range("H31").Value = 2698.5
$ws.Range("I31").Value = 1151
$ws.Range("J31").Value = 5019.75
$ws.Range("K31").Value = 1151
$ws.Range("L31").Value = 5019.75
$ws.Range("M31").Value = -856
$ws.Range("N31").Value = -5609.75
# row 34
$ws.Range("H34").Value = 2698.5
$ws.Range("I34").Value = 1151
$ws.Range("J34").Value = 5019.75
$ws.Range("K34").Value = 1151
$ws.Range("L34").Value = 5019.75
$ws.Range("M34").Value = -949
$ws.Range("N34").Value = -5423.75
# row 62
$ws.Range("H62").Value = 27514.25
$ws.Range("I62").Value = 2066.4285
$ws.Range("K62").Value = 2066.4285
$ws.Range("M62").Value = -1442.4285
# row 65
$ws.Range("H65").Value = 27514.25
$ws.Range("I65").Value = 2066.4285
$ws.Range("K65").Value = 10332.1425
$ws.Range("M65").Value = -7212.1425
# row 113
$ws.Range("H113").Value = 2023.6666
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 2023.6666
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 2023.6666
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -6363.6666
# row 132
$ws.Range("H132").Value = 9810746
$ws.Range("I132").Value = 11497971
$ws.Range("J132").Value = 24839.8
$ws.Range("K132").Value = 34493913
$ws.Range("L132").Value = 74519.39999999999
$ws.Range("M132").Value = -34491383
$ws.Range("N132").Value = -79579.39999999999

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# row 56
$ws.Range("H56").Value = 6815.609
$ws.Range("I56").Value = 6815.609
$ws.Range("K56").Value = 6815.609
$ws.Range("M56").Value = -6285.609
# row 86
$ws.Range("H86").Value = 205.75
$ws.Range("I86").Value = 195
$ws.Range("K86").Value = 585
$ws.Range("M86").Value = 601
# row 89
$ws.Range("H89").Value = 205.75
$ws.Range("I89").Value = 195
$ws.Range("K89").Value = 1755
$ws.Range("M89").Value = 4173
# row 132
$ws.Range("H132").Value = 2286.7273
$ws.Range("I132").Value = 1717
$ws.Range("J132").Value = 2761.5
$ws.Range("K132").Value = 15453
$ws.Range("L132").Value = 24853.5
$ws.Range("M132").Value = -12923
$ws.Range("N132").Value = -29913.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# row 113
$ws.Range("H113").Value = 13339.4
$ws.Range("I113").Value = 10000
$ws.Range("J113").Value = 14174.25
$ws.Range("K113").Value = 10000
$ws.Range("L113").Value = 14174.25
$ws.Range("M113").Value = -7830
$ws.Range("N113").Value = -18514.25
# row 122
$ws.Range("H122").Value = 848721.9
$ws.Range("I122").Value = 1574399.4
$ws.Range("K122").Value = 4723198.199999999
$ws.Range("M122").Value = -4720748.199999999
# row 126
$ws.Range("H126").Value = 3798.2334
$ws.Range("I126").Value = 2248.4546
$ws.Range("K126").Value = 6745.3638
$ws.Range("M126").Value = -4275.3638
# row 140
$ws.Range("H140").Value = 68689
$ws.Range("J140").Value = 68689
$ws.Range("L140").Value = 68689
$ws.Range("N140").Value = -79049

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# row 46
$ws.Range("H46").Value = 6284.6665
$ws.Range("J46").Value = 7174.905
$ws.Range("L46").Value = 7174.905
$ws.Range("N46").Value = -7550.905
# row 136
$ws.Range("H136").Value = 6126.04
$ws.Range("J136").Value = 7306.1333
$ws.Range("L136").Value = 21918.3999
$ws.Range("N136").Value = -27018.3999
# row 139
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# row 81
$ws.Range("H81").Value = 2203440.8
$ws.Range("I81").Value = 1897847.9
$ws.Range("K81").Value = 3795695.8
$ws.Range("M81").Value = -3794634.8
# row 84
$ws.Range("H84").Value = 2203440.8
$ws.Range("I84").Value = 1897847.9
$ws.Range("K84").Value = 18978479
$ws.Range("M84").Value = -18973175
